$d = $word.ActiveDocument

$replacements = @(
    @{old = "109×4=436";  new = "493×3=1479"},
    @{old = "553×3=1659"; new = "209×3=627"},
    @{old = "453×9=4077"; new = "994×3=2982"},
    @{old = "123×4=492";  new = "945×7=6615"},
    @{old = "730×3=2190"; new = "809×2=1618"},
    @{old = "507×3=1521"; new = "690×9=6210"},
    @{old = "394×2=788";  new = "757×3=2271"},
    @{old = "895×5=4475"; new = "533×5=2665"},
    @{old = "899×6=5394"; new = "443×2=886"},
    @{old = "680×3=2040"; new = "684×6=4104"},
    @{old = "633×4=2532"; new = "829×4=3316"},
    @{old = "599×6=3594"; new = "129×2=258"},
    @{old = "784×6=4704"; new = "817×5=4085"},
    @{old = "347×3=1041"; new = "926×9=8334"},
    @{old = "232×6=1392"; new = "607×4=2428"},
    @{old = "229×4=916";  new = "614×4=2456"},
    @{old = "337×2=674";  new = "478×5=2390"},
    @{old = "592×9=5328"; new = "475×2=950"},
    @{old = "816×5=4080"; new = "519×7=3633"},
    @{old = "721×9=6489"; new = "665×4=2660"},
    @{old = "533×4=2132"; new = "837×9=7533"},
    @{old = "563×2=1126"; new = "842×5=4210"},
    @{old = "498×8=3984"; new = "960×4=3840"},
    @{old = "148×8=1184"; new = "278×2=556"},
    @{old = "130×2=260";  new = "952×8=7616"}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
